$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(18, 8).Value = 14572.875
$ws.Cells.Item(18, 9).Value = 17916.834
$ws.Cells.Item(18, 11).Value = 17916.834
$ws.Cells.Item(18, 13).Value = -17632.834
$ws.Cells.Item(19, 8).Value = 3877.6155
$ws.Cells.Item(19, 9).Value = 6076.25
$ws.Cells.Item(19, 10).Value = 359.8
$ws.Cells.Item(19, 11).Value = 6076.25
$ws.Cells.Item(19, 12).Value = 359.8
$ws.Cells.Item(19, 13).Value = -5901.25
$ws.Cells.Item(19, 14).Value = -709.8
$ws.Cells.Item(40, 8).Value = 2775
$ws.Cells.Item(40, 9).Value = 1843.75
$ws.Cells.Item(40, 11).Value = 1843.75
$ws.Cells.Item(40, 13).Value = -1668.75
$ws.Cells.Item(58, 8).Value = 1419.375
$ws.Cells.Item(58, 9).Value = 313.75
$ws.Cells.Item(58, 10).Value = 2525
$ws.Cells.Item(58, 11).Value = 941.25
$ws.Cells.Item(58, 12).Value = 7575
$ws.Cells.Item(58, 13).Value = -791.25
$ws.Cells.Item(58, 14).Value = -7875
$ws.Cells.Item(70, 8).Value = 2835.7856
$ws.Cells.Item(70, 9).Value = 3138.7144
$ws.Cells.Item(70, 10).Value = 2532.8572
$ws.Cells.Item(70, 11).Value = 9416.143199999999
$ws.Cells.Item(70, 12).Value = 7598.571599999999
$ws.Cells.Item(70, 13).Value = -9146.143199999999
$ws.Cells.Item(70, 14).Value = -8138.571599999999
$ws.Cells.Item(73, 8).Value = 2835.7856
$ws.Cells.Item(73, 9).Value = 3138.7144
$ws.Cells.Item(73, 10).Value = 2532.8572
$ws.Cells.Item(73, 11).Value = 9416.143199999999
$ws.Cells.Item(73, 12).Value = 7598.571599999999
$ws.Cells.Item(73, 13).Value = -8480.143199999999
$ws.Cells.Item(73, 14).Value = -9470.571599999999
$ws.Cells.Item(82, 8).Value = 799.4
$ws.Cells.Item(82, 9).Value = 799.4
$ws.Cells.Item(82, 11).Value = 2398.2
$ws.Cells.Item(82, 13).Value = -1992.2
$ws.Cells.Item(85, 8).Value = 799.4
$ws.Cells.Item(85, 9).Value = 799.4
$ws.Cells.Item(85, 11).Value = 2398.2
$ws.Cells.Item(85, 13).Value = -994.1999999999998
$ws.Cells.Item(137, 8).Value = 1166.3125
$ws.Cells.Item(137, 9).Value = 1121.6
$ws.Cells.Item(137, 11).Value = 3364.8
$ws.Cells.Item(137, 13).Value = -814.7999999999997

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 1956.3529
$ws.Cells.Item(61, 9).Value = 1853.25
$ws.Cells.Item(61, 11).Value = 1853.25
$ws.Cells.Item(61, 13).Value = -1641.25
$ws.Cells.Item(74, 8).Value = 1729.4445
$ws.Cells.Item(74, 9).Value = 1684.1177
$ws.Cells.Item(74, 11).Value = 1684.1177
$ws.Cells.Item(74, 13).Value = -810.1177
$ws.Cells.Item(77, 8).Value = 1729.4445
$ws.Cells.Item(77, 9).Value = 1684.1177
$ws.Cells.Item(77, 11).Value = 8420.5885
$ws.Cells.Item(77, 13).Value = -4052.5885
$ws.Cells.Item(97, 8).Value = 2558.0476
$ws.Cells.Item(97, 9).Value = 2248.7368
$ws.Cells.Item(97, 11).Value = 2248.7368
$ws.Cells.Item(97, 13).Value = -1752.7368
$ws.Cells.Item(136, 8).Value = 1956.3529
$ws.Cells.Item(136, 9).Value = 1853.25
$ws.Cells.Item(136, 11).Value = 5559.75
$ws.Cells.Item(136, 13).Value = -3009.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 12003.9375
$ws.Cells.Item(20, 9).Value = 11760.4
$ws.Cells.Item(20, 10).Value = 12873.714
$ws.Cells.Item(20, 11).Value = 11760.4
$ws.Cells.Item(20, 12).Value = 12873.714
$ws.Cells.Item(20, 13).Value = -11513.4
$ws.Cells.Item(20, 14).Value = -13367.714
$ws.Cells.Item(134, 8).Value = 1405.1
$ws.Cells.Item(134, 9).Value = 1177.4584
$ws.Cells.Item(134, 11).Value = 3532.3752
$ws.Cells.Item(134, 13).Value = -997.3751999999999
$ws.Cells.Item(141, 8).Value = 60155.6
$ws.Cells.Item(141, 10).Value = 60155.6
$ws.Cells.Item(141, 12).Value = 60155.6
$ws.Cells.Item(141, 14).Value = -70515.60000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(3, 8).Value = 10000214
$ws.Cells.Item(3, 9).Value = 10000214
$ws.Cells.Item(3, 11).Value = 10000214
$ws.Cells.Item(3, 13).Value = -10000101
$ws.Cells.Item(7, 8).Value = 58827180
$ws.Cells.Item(7, 9).Value = 90911850
$ws.Cells.Item(7, 11).Value = 90911850
$ws.Cells.Item(7, 13).Value = -90911737
$ws.Cells.Item(31, 8).Value = 9980.796
$ws.Cells.Item(31, 9).Value = 2929.1428
$ws.Cells.Item(31, 11).Value = 2929.1428
$ws.Cells.Item(31, 13).Value = -2634.1428
$ws.Cells.Item(34, 8).Value = 9980.796
$ws.Cells.Item(34, 9).Value = 2929.1428
$ws.Cells.Item(34, 11).Value = 2929.1428
$ws.Cells.Item(34, 13).Value = -2727.1428
$ws.Cells.Item(58, 8).Value = 1455.6
$ws.Cells.Item(58, 9).Value = 1193
$ws.Cells.Item(58, 10).Value = 1849.5
$ws.Cells.Item(58, 11).Value = 1193
$ws.Cells.Item(58, 12).Value = 1849.5
$ws.Cells.Item(58, 13).Value = -990
$ws.Cells.Item(58, 14).Value = -2255.5
$ws.Cells.Item(62, 8).Value = 3979.8
$ws.Cells.Item(62, 9).Value = 4299.6665
$ws.Cells.Item(62, 11).Value = 4299.6665
$ws.Cells.Item(62, 13).Value = -3675.6665
$ws.Cells.Item(65, 8).Value = 3979.8
$ws.Cells.Item(65, 9).Value = 4299.6665
$ws.Cells.Item(65, 11).Value = 21498.3325
$ws.Cells.Item(65, 13).Value = -18378.3325
$ws.Cells.Item(68, 8).Value = 24999.555
$ws.Cells.Item(68, 10).Value = 24999.555
$ws.Cells.Item(68, 12).Value = 24999.555
$ws.Cells.Item(68, 14).Value = -26497.555
$ws.Cells.Item(71, 8).Value = 24999.555
$ws.Cells.Item(71, 10).Value = 24999.555
$ws.Cells.Item(71, 12).Value = 74998.66500000001
$ws.Cells.Item(71, 14).Value = -82486.66500000001
$ws.Cells.Item(136, 8).Value = 1455.6
$ws.Cells.Item(136, 9).Value = 1193
$ws.Cells.Item(136, 10).Value = 1849.5
$ws.Cells.Item(136, 11).Value = 3579
$ws.Cells.Item(136, 12).Value = 5548.5
$ws.Cells.Item(136, 13).Value = -1029
$ws.Cells.Item(136, 14).Value = -10648.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(10, 8).Value = 198.22223
$ws.Cells.Item(10, 9).Value = 160.5
$ws.Cells.Item(10, 10).Value = 500
$ws.Cells.Item(10, 11).Value = 481.5
$ws.Cells.Item(10, 12).Value = 1500
$ws.Cells.Item(10, 13).Value = -342.5
$ws.Cells.Item(10, 14).Value = -1778
$ws.Cells.Item(19, 8).Value = 55
$ws.Cells.Item(19, 9).Value = 10
$ws.Cells.Item(19, 10).Value = 100
$ws.Cells.Item(19, 11).Value = 30
$ws.Cells.Item(19, 12).Value = 300
$ws.Cells.Item(19, 13).Value = 144
$ws.Cells.Item(19, 14).Value = -648
$ws.Cells.Item(23, 8).Value = 112.583336
$ws.Cells.Item(23, 9).Value = 118.166664
$ws.Cells.Item(23, 10).Value = 107
$ws.Cells.Item(23, 11).Value = 354.499992
$ws.Cells.Item(23, 12).Value = 321
$ws.Cells.Item(23, 13).Value = -119.499992
$ws.Cells.Item(23, 14).Value = -791
$ws.Cells.Item(34, 8).Value = 875.38464
$ws.Cells.Item(34, 9).Value = 193.16667
$ws.Cells.Item(34, 11).Value = 579.50001
$ws.Cells.Item(34, 13).Value = -495.50001
$ws.Cells.Item(41, 8).Value = 199.27586
$ws.Cells.Item(41, 9).Value = 109.9375
$ws.Cells.Item(41, 10).Value = 309.23077
$ws.Cells.Item(41, 11).Value = 329.8125
$ws.Cells.Item(41, 12).Value = 927.69231
$ws.Cells.Item(41, 13).Value = 8.1875
$ws.Cells.Item(41, 14).Value = -1603.69231
$ws.Cells.Item(81, 8).Value = 35723410
$ws.Cells.Item(81, 9).Value = 7217
$ws.Cells.Item(81, 11).Value = 21651
$ws.Cells.Item(81, 13).Value = -20528
$ws.Cells.Item(84, 8).Value = 35723410
$ws.Cells.Item(84, 9).Value = 7217
$ws.Cells.Item(84, 11).Value = 64953
$ws.Cells.Item(84, 13).Value = -59337
$ws.Cells.Item(104, 8).Value = 4187.5557
$ws.Cells.Item(104, 9).Value = 2700
$ws.Cells.Item(104, 10).Value = 4373.5
$ws.Cells.Item(104, 11).Value = 8100
$ws.Cells.Item(104, 12).Value = 13120.5
$ws.Cells.Item(104, 13).Value = -5479
$ws.Cells.Item(104, 14).Value = -18362.5
$ws.Cells.Item(140, 8).Value = 1162.2609
$ws.Cells.Item(140, 9).Value = 1124.409
$ws.Cells.Item(140, 11).Value = 3373.227
$ws.Cells.Item(140, 13).Value = 1806.773

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(40, 8).Value = 22900
$ws.Cells.Item(40, 9).Value = 22900
$ws.Cells.Item(40, 11).Value = 22900
$ws.Cells.Item(40, 13).Value = -22749
$ws.Cells.Item(97, 8).Value = 34212.54
$ws.Cells.Item(97, 9).Value = 22215.566
$ws.Cells.Item(97, 10).Value = 74202.44500000001
$ws.Cells.Item(97, 11).Value = 22215.566
$ws.Cells.Item(97, 12).Value = 74202.44500000001
$ws.Cells.Item(97, 13).Value = -21719.566
$ws.Cells.Item(97, 14).Value = -75194.44500000001
$ws.Cells.Item(126, 8).Value = 1738.75
$ws.Cells.Item(126, 9).Value = 1489
$ws.Cells.Item(126, 10).Value = 1988.5
$ws.Cells.Item(126, 11).Value = 4467
$ws.Cells.Item(126, 12).Value = 5965.5
$ws.Cells.Item(126, 13).Value = -1997
$ws.Cells.Item(126, 14).Value = -10905.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(43, 8).Value = 10077492
$ws.Cells.Item(43, 10).Value = 10077492
$ws.Cells.Item(43, 12).Value = 10077492
$ws.Cells.Item(43, 14).Value = -10077878
$ws.Cells.Item(46, 8).Value = 21646.592
$ws.Cells.Item(46, 9).Value = 44252.3
$ws.Cells.Item(46, 10).Value = 2808.5
$ws.Cells.Item(46, 11).Value = 44252.3
$ws.Cells.Item(46, 12).Value = 2808.5
$ws.Cells.Item(46, 13).Value = -44064.3
$ws.Cells.Item(46, 14).Value = -3184.5
$ws.Cells.Item(68, 8).Value = 3489.111
$ws.Cells.Item(68, 9).Value = 3507
$ws.Cells.Item(68, 11).Value = 3507
$ws.Cells.Item(68, 13).Value = -2758
$ws.Cells.Item(71, 8).Value = 3489.111
$ws.Cells.Item(71, 9).Value = 3507
$ws.Cells.Item(71, 11).Value = 17535
$ws.Cells.Item(71, 13).Value = -13791
$ws.Cells.Item(107, 8).Value = 3147.5
$ws.Cells.Item(107, 9).Value = 3147.5
$ws.Cells.Item(107, 11).Value = 3147.5
$ws.Cells.Item(107, 13).Value = -1227.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(41, 8).Value = 15203.25
$ws.Cells.Item(41, 9).Value = 0
$ws.Cells.Item(41, 10).Value = 15203.25
$ws.Cells.Item(41, 11).Value = 0
$ws.Cells.Item(41, 12).Value = 15203.25
$ws.Cells.Item(41, 13).ClearContents()
$ws.Cells.Item(41, 14).Value = -15983.25
